# Update the "dSF" (column F) values for the lorenzen_michael 2024 data sheet.
# This mirrors the author's "repull data, push all data, mean calculation" commit,
# which recalculated the dSF figures while leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = 0
    5  = 1
    6  = 2
    7  = 3
    8  = 3
    9  = 2
    10 = 1
    11 = -4
    12 = -3
    13 = -1
    14 = -3
    15 = 3
    16 = -3
    17 = 2
    18 = 3
    20 = -2
    21 = 5
    22 = 2
    23 = 1
    25 = -2
    27 = -3
    28 = 1
    29 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
